{"js": "// The diff collapses three paragraphs (Title, Author, Abstract) whose text\n// had been split word-by-word across many runs into a single run each,\n// leaving the rendered text itself unchanged. Re-set each paragraph's text\n// in place (Range.insertText with \"Replace\") so Word coalesces all of the\n// runs in that paragraph into one run carrying the full sentence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nconst replacements = {\n  \"Title\": \"Answers: Introduction to vectors\",\n  \"Author\": \"Zheng Chen\",\n  \"Abstract\": \"Answers to questions relating to the guide on introduction to vectors.\",\n};\n\nfor (const paragraph of paragraphs.items) {\n  const newText = replacements[paragraph.style];\n  if (newText !== undefined) {\n    paragraph.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The diff collapses three paragraphs (Title, Author, Abstract) whose text\n# had been split word-by-word across many runs into a single run each; the\n# rendered text itself is unchanged. Re-running Find/Replace scoped to each\n# paragraph's own Range (not the whole document) makes Word coalesce all of\n# the runs in that paragraph into a single run carrying the full sentence,\n# without touching any other paragraph that might contain similar text\n# (e.g. \"Zheng Chen\" also appears, already as a single run, further down).\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($paragraph, [string]$text) {\n    $range = $paragraph.Range\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $text\n    $find.Replacement.Text = $text\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nforeach ($paragraph in $d.Paragraphs) {\n    $styleName = $paragraph.Style.NameLocal\n    if ($styleName -eq \"Title\") {\n        Set-ParagraphText $paragraph \"Answers: Introduction to vectors\"\n    } elseif ($styleName -eq \"Author\") {\n        Set-ParagraphText $paragraph \"Zheng Chen\"\n    } elseif ($styleName -eq \"Abstract\") {\n        Set-ParagraphText $paragraph \"Answers to questions relating to the guide on introduction to vectors.\"\n    }\n}\n"}
